$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.673.89"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.626.64"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'212.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'22.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").Value = "'0.0611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.857.20"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.631.38"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "'0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "'64.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "27.683.94"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "'230.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'7.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "'4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'9.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("D25").Value = "'149.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'6.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D31").Value = "'0.0483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "1.461.99"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'0.870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'69.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.19%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "1.767.68"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").Value = "'1.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "'86.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
